# This workbook contains a weekly "Ajo" (garlic) price dataset for
# "Terminal La Palmera de La Serena". The commit adds one new week's worth
# of data (2 rows) at the top of the historical block that starts at row 220,
# pushing all the existing rows in that block down by 2 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 220-221; this shifts old rows 220:338 down to 222:340
# and keeps everything else (rows 1:219) untouched.
$ws.Range("A220:A221").EntireRow.Insert()

# Populate the first new row (220) with the newest price record.
$ws.Cells.Item(220, 1).Value  = 8
$ws.Cells.Item(220, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(220, 3).Value  = "Coquimbo"
$ws.Cells.Item(220, 4).Value  = 44830
$ws.Cells.Item(220, 5).Value  = 4
$ws.Cells.Item(220, 6).Value  = 100112003
$ws.Cells.Item(220, 7).Value  = "Ajo"
$ws.Cells.Item(220, 8).Value  = "Chino"
$ws.Cells.Item(220, 9).Value  = "Primera"
$ws.Cells.Item(220, 10).Value = 500
$ws.Cells.Item(220, 11).Value = 22000
$ws.Cells.Item(220, 12).Value = 23000
$ws.Cells.Item(220, 13).Value = 22500
$ws.Cells.Item(220, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(220, 15).Value = "China"
$ws.Cells.Item(220, 16).Value = 2250
$ws.Cells.Item(220, 17).Value = 10
$ws.Cells.Item(220, 18).Value = "Hortaliza"

# Populate the second new row (221) with the second newest price record.
$ws.Cells.Item(221, 1).Value  = 8
$ws.Cells.Item(221, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(221, 3).Value  = "Coquimbo"
$ws.Cells.Item(221, 4).Value  = 44830
$ws.Cells.Item(221, 5).Value  = 4
$ws.Cells.Item(221, 6).Value  = 100112003
$ws.Cells.Item(221, 7).Value  = "Ajo"
$ws.Cells.Item(221, 8).Value  = "Chino"
$ws.Cells.Item(221, 9).Value  = "Primera"
$ws.Cells.Item(221, 10).Value = 400
$ws.Cells.Item(221, 11).Value = 24000
$ws.Cells.Item(221, 12).Value = 25000
$ws.Cells.Item(221, 13).Value = 24500
$ws.Cells.Item(221, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(221, 15).Value = "China"
$ws.Cells.Item(221, 16).Value = 2450
$ws.Cells.Item(221, 17).Value = 10
$ws.Cells.Item(221, 18).Value = "Hortaliza"

# Make sure the D column keeps the date/time number format used by the rest
# of the column (it is normally inherited from the row that was duplicated
# by Insert(), but set it explicitly to be safe).
$ws.Range("D220:D221").NumberFormat = $ws.Range("D222").NumberFormat
